$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 1.168007333333333
$ws.Range("N2").Value = 3.504022
$ws.Range("O2").Value = 0.1638609704511517
$ws.Range("P2").Value = 0.1638609704511517
$ws.Range("Q2").Value = 41.74324044224312
$ws.Range("R2").Value = 375.689163980188
$ws.Range("S2").Value = 0.003194845530630382
$ws.Range("T2").Value = 0.003194845530630382

$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.5019752511630595
$ws.Range("P3").Value = 0.5019752511630595
$ws.Range("Q3").Value = 127.8771482169487
$ws.Range("R3").Value = 1150.894333952538
$ws.Range("S3").Value = 0.009787159097434068
$ws.Range("T3").Value = 0.009787159097434066

$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 1.915392333333333
$ws.Range("N4").Value = 5.746177
$ws.Range("O4").Value = 0.2687123938160456
$ws.Range("P4").Value = 0.2687123938160456
$ws.Range("Q4").Value = 68.45392184600645
$ws.Range("R4").Value = 616.085296614058
$ws.Range("S4").Value = 0.005239164567648576
$ws.Range("T4").Value = 0.005239164567648576

$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.46654
$ws.Range("N5").Value = 1.39962
$ws.Range("O5").Value = 0.06545138456974327
$ws.Range("P5").Value = 0.06545138456974327
$ws.Range("Q5").Value = 16.67360370105333
$ws.Range("R5").Value = 150.06243330948
$ws.Range("S5").Value = 0.001276124893502636
$ws.Range("T5").Value = 0.001276124893502636

$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 1.168007333333333
$ws.Range("N6").Value = 3.504022
$ws.Range("O6").Value = 0.1638609704511517
$ws.Range("P6").Value = 0.1638609704511517
$ws.Range("Q6").Value = 1973.102298308249
$ws.Range("R6").Value = 17757.92068477424
$ws.Range("S6").Value = 0.1510126428241397
$ws.Range("T6").Value = 0.1510126428241397

$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.5019752511630595
$ws.Range("P7").Value = 0.5019752511630595
$ws.Range("Q7").Value = 6044.444378894692
$ws.Range("R7").Value = 54399.99941005222
$ws.Range("S7").Value = 0.4626154056193809
$ws.Range("T7").Value = 0.4626154056193809

$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 1.915392333333333
$ws.Range("N8").Value = 5.746177
$ws.Range("O8").Value = 0.2687123938160456
$ws.Range("P8").Value = 0.2687123938160456
$ws.Range("Q8").Value = 3235.651786771316
$ws.Range("R8").Value = 29120.86608094184
$ws.Range("S8").Value = 0.2476426731639489
$ws.Range("T8").Value = 0.247642673163949

$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.46654
$ws.Range("N9").Value = 1.39962
$ws.Range("O9").Value = 0.06545138456974327
$ws.Range("P9").Value = 0.06545138456974327
$ws.Range("Q9").Value = 788.1210331322666
$ws.Range("R9").Value = 7093.089298190399
$ws.Range("S9").Value = 0.06031934592577399
$ws.Range("T9").Value = 0.06031934592577399

$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 1.168007333333333
$ws.Range("N10").Value = 3.504022
$ws.Range("O10").Value = 0.1638609704511517
$ws.Range("P10").Value = 0.1638609704511517
$ws.Range("Q10").Value = 109.374020608722
$ws.Range("R10").Value = 984.3661854784981
$ws.Range("S10").Value = 0.008371010424845534
$ws.Range("T10").Value = 0.008371010424845534

$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.5019752511630595
$ws.Range("P11").Value = 0.5019752511630595
$ws.Range("Q11").Value = 335.0587471477469
$ws.Range("R11").Value = 3015.528724329723
$ws.Range("S11").Value = 0.02564393490976601
$ws.Range("T11").Value = 0.02564393490976601

$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 1.915392333333333
$ws.Range("N12").Value = 5.746177
$ws.Range("O12").Value = 0.2687123938160456
$ws.Range("P12").Value = 0.2687123938160456
$ws.Range("Q12").Value = 179.360312697627
$ws.Range("R12").Value = 1614.242814278643
$ws.Range("S12").Value = 0.01372745592636337
$ws.Range("T12").Value = 0.01372745592636337

$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.46654
$ws.Range("N13").Value = 1.39962
$ws.Range("O13").Value = 0.06545138456974327
$ws.Range("P13").Value = 0.06545138456974327
$ws.Range("Q13").Value = 43.68753013662
$ws.Range("R13").Value = 393.18777122958
$ws.Range("S13").Value = 0.00334365298243627
$ws.Range("T13").Value = 0.00334365298243627

$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 1.168007333333333
$ws.Range("N14").Value = 3.504022
$ws.Range("O14").Value = 0.1638609704511517
$ws.Range("P14").Value = 0.1638609704511517
$ws.Range("Q14").Value = 16.75652948852689
$ws.Range("R14").Value = 150.808765396742
$ws.Range("S14").Value = 0.001282471671536087
$ws.Range("T14").Value = 0.001282471671536087

$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.5019752511630595
$ws.Range("P15").Value = 0.5019752511630595
$ws.Range("Q15").Value = 51.33231589844633
$ws.Range("R15").Value = 461.9908430860169
$ws.Range("S15").Value = 0.003928751536478597
$ws.Range("T15").Value = 0.003928751536478597

$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 1.915392333333333
$ws.Range("N16").Value = 5.746177
$ws.Range("O16").Value = 0.2687123938160456
$ws.Range("P16").Value = 0.2687123938160456
$ws.Range("Q16").Value = 27.47870428518856
$ws.Range("R16").Value = 247.308338566697
$ws.Range("S16").Value = 0.002103100158084687
$ws.Range("T16").Value = 0.002103100158084687

$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.46654
$ws.Range("N17").Value = 1.39962
$ws.Range("O17").Value = 0.06545138456974327
$ws.Range("P17").Value = 0.06545138456974327
$ws.Range("Q17").Value = 6.693101185646666
$ws.Range("R17").Value = 60.23791067081999
$ws.Range("S17").Value = 0.0005122607680303773
$ws.Range("T17").Value = 0.0005122607680303773
